# Generate Report for Handback
#
# This mirrors the localization tool's "handback" run: the two files that
# were previously only "handed off" (f21df298...md and fa6bbf84...md) now
# come back from de-de fully in sync, and the zh-cn/de-de per-file report
# rows get their "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns populated for the first time.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    everywhere it appears (Overview!E:F, zh-cn!C, de-de!C).
# ---------------------------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn: fill in "Latest Target File" (I) / "Latest Handback File" (J)
#    for both rows; the handback datetime (K) text itself is corrected
#    (it had the zero-date placeholder).
# ---------------------------------------------------------------------
$zhcn.Range("I2").Value = "f21df298-b515-44a4-8844-3ba8154e21d3.md"
$zhcn.Range("J2").Value = "f21df298-b515-44a4-8844-3ba8154e21d3.981de31c99c112ddc50525003931eb553741234b.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-10-14 08:20:23"

$zhcn.Range("I3").Value = "fa6bbf84-bb9c-4279-bdc9-f7b1323899a5.md"
$zhcn.Range("J3").Value = "fa6bbf84-bb9c-4279-bdc9-f7b1323899a5.bce96e7820086f195e65185f0295395ff792c021.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-10-14 08:20:23"

# ---------------------------------------------------------------------
# 3. de-de: same three columns, but this locale actually completed a
#    handback, so it gets the real handback timestamp.
# ---------------------------------------------------------------------
$dede.Range("I2").Value = "f21df298-b515-44a4-8844-3ba8154e21d3.md"
$dede.Range("J2").Value = "f21df298-b515-44a4-8844-3ba8154e21d3.981de31c99c112ddc50525003931eb553741234b.de-de.xlf"
$dede.Range("K2").Value = "2016-10-14 08:20:39"

$dede.Range("I3").Value = "fa6bbf84-bb9c-4279-bdc9-f7b1323899a5.md"
$dede.Range("J3").Value = "fa6bbf84-bb9c-4279-bdc9-f7b1323899a5.bce96e7820086f195e65185f0295395ff792c021.de-de.xlf"
$dede.Range("K3").Value = "2016-10-14 08:20:39"

# ---------------------------------------------------------------------
# 4. Hyperlinks: "Latest Target File" points at the same source .md page
#    as the "Source File Name" column (A). Rebuild each sheet's
#    hyperlinks in row order so the new I2/I3 links land right after
#    their row's existing A-column link (matching how the report
#    generator laid out the relationships).
# ---------------------------------------------------------------------
function Set-HandbackHyperlinks($ws) {
    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add(
        $ws.Range("A2"),
        "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/39efdfeae1c0ebd5954af908d9c7cf6c1df6afcd/e2e/f21df298-b515-44a4-8844-3ba8154e21d3.md",
        [Type]::Missing, [Type]::Missing,
        "f21df298-b515-44a4-8844-3ba8154e21d3.md")

    $ws.Hyperlinks.Add(
        $ws.Range("I2"),
        "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/39efdfeae1c0ebd5954af908d9c7cf6c1df6afcd/e2e/f21df298-b515-44a4-8844-3ba8154e21d3.md",
        [Type]::Missing, [Type]::Missing,
        "f21df298-b515-44a4-8844-3ba8154e21d3.md")

    $ws.Hyperlinks.Add(
        $ws.Range("A3"),
        "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/39efdfeae1c0ebd5954af908d9c7cf6c1df6afcd/e2e/fa6bbf84-bb9c-4279-bdc9-f7b1323899a5.md",
        [Type]::Missing, [Type]::Missing,
        "fa6bbf84-bb9c-4279-bdc9-f7b1323899a5.md")

    $ws.Hyperlinks.Add(
        $ws.Range("I3"),
        "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/39efdfeae1c0ebd5954af908d9c7cf6c1df6afcd/e2e/fa6bbf84-bb9c-4279-bdc9-f7b1323899a5.md",
        [Type]::Missing, [Type]::Missing,
        "fa6bbf84-bb9c-4279-bdc9-f7b1323899a5.md")
}

Set-HandbackHyperlinks $zhcn
Set-HandbackHyperlinks $dede

# ---------------------------------------------------------------------
# 5. Column widths widened to fit the newly-populated file-name /
#    datetime columns.
# ---------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 29.1
$overview.Columns.Item(6).ColumnWidth = 29.1

foreach ($ws in @($zhcn, $dede)) {
    $ws.Columns.Item(3).ColumnWidth  = 29.1
    $ws.Columns.Item(9).ColumnWidth  = 39.1
    $ws.Columns.Item(10).ColumnWidth = 39.1
}
